# New crime data collected - weekly CompStat refresh (81st Precinct)
# Updates the report header (volume/week-of dates), refreshes the weekly
# crime-stat grid (rows 14-30) with newly collected figures, and inserts a
# blank spacer row above the "Prepared by" footer block (old rows 56/57
# shift down to 57/58).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header: volume number + week-covering dates
# ---------------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 31   Number  51"
$ws.Range("C9").Value = "Report Covering the Week  12/16/2024  Through  12/22/2024"

# ---------------------------------------------------------------------------
# Crime complaint grid refresh (rows 14-30)
# ---------------------------------------------------------------------------

# Row 14 - Murder
$ws.Range("G14").Value = 1

# Row 15 - Rape
$ws.Range("C15").Value = 0
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = "***.*"
$ws.Range("F15").Value = 2
$ws.Range("H15").Value = 100

# Row 16 - Robbery
$ws.Range("C16").Value = 2
$ws.Range("E16").Value = -33.333333333333
$ws.Range("F16").Value = 9
$ws.Range("G16").Value = 13
$ws.Range("H16").Value = -30.76923076923
$ws.Range("I16").Value = 127
$ws.Range("J16").Value = 150
$ws.Range("K16").Value = -15.333333333333
$ws.Range("L16").Value = -16.447368421052
$ws.Range("M16").Value = -63.400576368876
$ws.Range("N16").Value = -89.451827242524

# Row 17 - Fel. Assault
$ws.Range("D17").Value = 8
$ws.Range("E17").Value = -87.5
$ws.Range("F17").Value = 14
$ws.Range("G17").Value = 28
$ws.Range("H17").Value = -50
$ws.Range("J17").Value = 299
$ws.Range("K17").Value = -6.354515050167
$ws.Range("L17").Value = -14.110429447852
$ws.Range("M17").Value = -12.225705329153
$ws.Range("N17").Value = -66.94214876033

# Row 18 - Burglary
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = 0
$ws.Range("G18").Value = 8
$ws.Range("H18").Value = -25
$ws.Range("I18").Value = 79
$ws.Range("J18").Value = 135
$ws.Range("K18").Value = -41.481481481481
$ws.Range("L18").Value = -60.5
$ws.Range("M18").Value = -67.622950819672
$ws.Range("N18").Value = -90.628706998813

# Row 19 - Gr. Larceny
$ws.Range("C19").Value = 3
$ws.Range("D19").Value = 2
$ws.Range("E19").Value = 50
$ws.Range("F19").Value = 24
$ws.Range("G19").Value = 25
$ws.Range("H19").Value = -4
$ws.Range("I19").Value = 268
$ws.Range("J19").Value = 336
$ws.Range("K19").Value = -20.238095238095
$ws.Range("L19").Value = -37.236533957845
$ws.Range("M19").Value = -25.348189415041
$ws.Range("N19").Value = -28.912466843501

# Row 20 - G.L.A.
$ws.Range("C20").Value = 0
$ws.Range("E20").Value = -100
$ws.Range("F20").Value = 1
$ws.Range("G20").Value = 8
$ws.Range("H20").Value = -87.5
$ws.Range("I20").Value = 87
$ws.Range("J20").Value = 106
$ws.Range("K20").Value = -17.924528301886
$ws.Range("L20").Value = -35.074626865671
$ws.Range("M20").Value = 6.097560975609
$ws.Range("N20").Value = -85.077186963979

# Row 21 - TOTAL
$ws.Range("C21").Value = 8
$ws.Range("D21").Value = 18
$ws.Range("E21").Value = -55.555555555555
$ws.Range("F21").Value = 56
$ws.Range("G21").Value = 84
$ws.Range("H21").Value = -33.333333333333
$ws.Range("I21").Value = 858
$ws.Range("J21").Value = 1050
$ws.Range("K21").Value = -18.285714285714
$ws.Range("L21").Value = -32.173913043478
$ws.Range("M21").Value = -38.362068965517
$ws.Range("N21").Value = -78.322385042951

# Row 22 - Transit
$ws.Range("G22").Value = 2
$ws.Range("H22").Value = 100
$ws.Range("M22").Value = -24.137931034482

# Row 23 - Housing
$ws.Range("C23").Value = 1
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 5
$ws.Range("G23").Value = 6
$ws.Range("H23").Value = -16.666666666666
$ws.Range("I23").Value = 86
$ws.Range("J23").Value = 85
$ws.Range("K23").Value = 1.176470588235
$ws.Range("L23").Value = 4.878048780487
$ws.Range("M23").Value = 2.380952380952

# Row 24 - Petit Larceny
$ws.Range("C24").Value = 21
$ws.Range("D24").Value = 18
$ws.Range("E24").Value = 16.666666666666
$ws.Range("F24").Value = 59
$ws.Range("G24").Value = 61
$ws.Range("H24").Value = -3.27868852459
$ws.Range("I24").Value = 719
$ws.Range("J24").Value = 827
$ws.Range("K24").Value = -13.059250302297
$ws.Range("L24").Value = -17.165898617511
$ws.Range("M24").Value = -8.524173027989

# Row 25 - Retail Theft
$ws.Range("C25").Value = 2
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = "***.*"
$ws.Range("G25").Value = 5
$ws.Range("H25").Value = 80
$ws.Range("I25").Value = 132
$ws.Range("K25").Value = -28.648648648648
$ws.Range("L25").Value = -48.638132295719

# Row 26 - Misd. Assault
$ws.Range("C26").Value = 11
$ws.Range("D26").Value = 7
$ws.Range("E26").Value = 57.142857142857
$ws.Range("F26").Value = 31
$ws.Range("G26").Value = 31
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 409
$ws.Range("J26").Value = 485
$ws.Range("K26").Value = -15.670103092783
$ws.Range("L26").Value = 0.738916256157
$ws.Range("M26").Value = -47.831632653061

# Row 27 - UCR Rape*
$ws.Range("C27").Value = 0
$ws.Range("D27").Value = 0
$ws.Range("E27").Value = "***.*"
$ws.Range("F27").Value = 2
$ws.Range("H27").Value = 100

# Row 28 - Other Sex Crimes
$ws.Range("C28").Value = 1
$ws.Range("F28").Value = 5
$ws.Range("H28").Value = 400
$ws.Range("L28").Value = 17.647058823529

# Row 29 - Shooting Vic.
$ws.Range("D29").Value = 0
$ws.Range("E29").Value = "***.*"
$ws.Range("G29").Value = 5
$ws.Range("N29").Value = -89.017341040462

# Row 30 - Shooting Inc.
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = "***.*"
$ws.Range("G30").Value = 3
$ws.Range("N30").Value = -89.743589743589

# ---------------------------------------------------------------------------
# Insert a blank spacer row above the footer block (old row 56 -> 57,
# old row 57 -> 58); dimension grows to A1:N58.
# ---------------------------------------------------------------------------
$ws.Rows.Item(56).Insert()
$ws.Cells.Item(56, 1).ClearFormats()
$ws.Cells.Item(56, 1).ClearContents()
